$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 191 (Vega Modelo de Temuco - Poroto verde
# dataset), pushing the current row 191 (and everything below it) down by one.
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new weekly price observation.
# (Mercado/Region/Codreg/Categoria/Variedad/Calidad/Clasificacion columns follow
# the same pattern as the surrounding rows of this subset.)
$ws.Range("A191").Value = 10
$ws.Range("B191").Value = "Vega Modelo de Temuco"
$ws.Range("C191").Value = "La Araucanía"
$ws.Range("D191").Value = 45211
$ws.Range("D191").NumberFormat = $ws.Range("D192").NumberFormat
$ws.Range("E191").Value = 9
$ws.Range("F191").Value = 100112031
$ws.Range("G191").Value = "Poroto verde"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 100
$ws.Range("K191").Value = 1600
$ws.Range("L191").Value = 1600
$ws.Range("M191").Value = 1600
$ws.Range("N191").Value = "$/kilo"
$ws.Range("O191").Value = "Perú"
$ws.Range("P191").Value = 1600
$ws.Range("Q191").Value = 1
$ws.Range("R191").Value = "Hortaliza"
